$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild the Team / PER (Player Efficiency Rating) table: the previous
# shared-string ordering misaligned team names with their stat row (the
# "PER bug"), and the PER values themselves were computed on the wrong
# scale. Re-write each row with the corrected team label and corrected
# PER value.

$ws.Range("B2").Value = "POR"
$ws.Range("C2").Value = 12.78461538461538

$ws.Range("B3").Value = "NJN"
$ws.Range("C3").Value = 12.29411764705883

$ws.Range("B4").Value = "CLE"
$ws.Range("C4").Value = 15.4

$ws.Range("B5").Value = "DAL"
$ws.Range("C5").Value = 11.44375

$ws.Range("B6").Value = "MIA"
$ws.Range("C6").Value = 13.53846153846154

$ws.Range("B7").Value = "SEA"
$ws.Range("C7").Value = 15.46363636363636

$ws.Range("B8").Value = "ATL"
$ws.Range("C8").Value = 12.87692307692308

$ws.Range("B9").Value = "MIL"
$ws.Range("C9").Value = 14.25384615384615

$ws.Range("B10").Value = "LAC"
$ws.Range("C10").Value = 12.65714285714286

$ws.Range("B11").Value = "DET"
$ws.Range("C11").Value = 12.94545454545455

$ws.Range("B12").Value = "SAS"
$ws.Range("C12").Value = 12.02857142857143

$ws.Range("B13").Value = "ORL"
$ws.Range("C13").Value = 14.08461538461538

$ws.Range("B14").Value = "UTA"
$ws.Range("C14").Value = 13.25384615384615

$ws.Range("B15").Value = "HOU"
$ws.Range("C15").Value = 12.93076923076923

$ws.Range("B16").Value = "DEN"
$ws.Range("C16").Value = 12.34615384615385

$ws.Range("B17").Value = "LAL"
$ws.Range("C17").Value = 13.55

$ws.Range("B18").Value = "GSW"
$ws.Range("C18").Value = 13.30588235294118

$ws.Range("B19").Value = "IND"
$ws.Range("C19").Value = 14.35384615384615

$ws.Range("B20").Value = "CHI"
$ws.Range("C20").Value = 13.79285714285714

$ws.Range("B21").Value = "PHI"
$ws.Range("C21").Value = 13.52142857142857

$ws.Range("B22").Value = "CHH"
$ws.Range("C22").Value = 14.06363636363636

$ws.Range("B23").Value = "BOS"
$ws.Range("C23").Value = 14.68571428571429

$ws.Range("B24").Value = "WSB"
$ws.Range("C24").Value = 13.26

$ws.Range("B25").Value = "SAC"
$ws.Range("C25").Value = 13.74615384615385

$ws.Range("B26").Value = "PHO"
$ws.Range("C26").Value = 14.85384615384615

$ws.Range("B27").Value = "NYK"
$ws.Range("C27").Value = 12.9

$ws.Range("B28").Value = "MIN"
$ws.Range("C28").Value = 12.51538461538462
